# "Don't need to convert with fuhep"
# The Clint (column N) formula previously converted Clint hep (column J)
# to uL/min/10^6 cells by multiplying by fuhep (column I):
#     =J{row}*I{row}/120/40*1000
# That multiplication by fuhep is no longer wanted, so the formula becomes:
#     =J{row}/120/40*1000
# Every row in the table (2-58) should show this corrected formula/value,
# including rows that previously had no computed Clint value at all.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 uses a standalone formula (not shared) in the source workbook.
$ws.Range("N2").Formula = "=J2/120/40*1000"

# Rows 3-58: set (or add) the corrected formula for every row, including
# rows that previously had no value in column N at all.
for ($r = 3; $r -le 58; $r++) {
    $ws.Cells.Item($r, 14).Formula = "=J$r/120/40*1000"
}

# Give the whole recalculated column a numeric display format (0.00),
# which creates/uses a new cell style.
$ws.Range("N2:N58").NumberFormat = "0.00"

# Update the sheet's active selection / view state.
$ws.Range("O8").Select()
